# Generate Report for Handoff
#
# A new handoff package (b.*.xlf) was produced for the "b.md" source file.
# It supersedes the previous "in sync with en-US" handback status, and the
# new handoff is flagged as based on a stale handback version, so an error
# detail message is recorded. Reflect this across the Overview summary
# sheet and both per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/2c527ef1d57486219f07e1851e208f337201f24b/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/3e4d35dfa3e2d1ae6444dc89df79a896acff6a58/e2e/b.md."

# ---- Overview sheet : row 3 is the "b.md" file ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $status
$overview.Range("F3").Value = $status
$overview.Range("G3").Value = "2016-08-14 00:49:30"

# ---- zh-cn sheet : row 3 is the "b.md" file ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $status
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-14 00:49:22"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = (40 - 5/6)

# ---- de-de sheet : row 3 is the "b.md" file ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $status
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-14 00:49:30"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = (40 - 5/6)
